$wb = $excel.ActiveWorkbook

# --- Remove stray empty inline-string cells in column B of "ODI Batting" ---
# (these were placeholder blanks written where INNING_NUMBER was not applicable;
#  clearing them drops the empty <c> elements entirely, matching upstream)
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$emptyBRows = @(4,5,8,10,11,14,16,17,18,19,20,21,23,25,32,34)
foreach ($r in $emptyBRows) {
    $wsBatting.Cells.Item($r, 2).Clear()
}

# --- Add the two new "extra" sheets at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "ODI Batting Extra"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "ODI Bowling Extra"

$data4 = @(
    ,@(1,1,'MATCH_CODE','txt')
    ,@(1,2,'BATTING_POSITION','txt')
    ,@(1,3,'NUM_4','txt')
    ,@(1,4,'NUM_6','txt')
    ,@(1,5,'PERCENT_RUNS_OF_TOTAL','txt')
    ,@(1,6,'MAN_OF_MATCH','txt')
    ,@(2,1,'2964','forcetxt')
    ,@(2,6,'NO','txt')
    ,@(3,1,'3007','forcetxt')
    ,@(3,6,'NO','txt')
    ,@(4,1,'3008','forcetxt')
    ,@(4,6,'NO','txt')
    ,@(5,1,'3135','forcetxt')
    ,@(5,2,'10','num')
    ,@(5,3,'1','forcetxt')
    ,@(5,4,'0','forcetxt')
    ,@(5,5,'8.59%','forcetxt')
    ,@(5,6,'NO','txt')
    ,@(6,1,'3138','forcetxt')
    ,@(6,2,'10','num')
    ,@(6,3,'0','forcetxt')
    ,@(6,4,'1','forcetxt')
    ,@(6,5,'2.56%','forcetxt')
    ,@(6,6,'NO','txt')
    ,@(7,1,'3142','forcetxt')
    ,@(7,6,'NO','txt')
    ,@(8,1,'3147','forcetxt')
    ,@(8,6,'NO','txt')
    ,@(9,1,'3149','forcetxt')
    ,@(9,2,'10','num')
    ,@(9,6,'NO','txt')
    ,@(10,1,'3153','forcetxt')
    ,@(10,2,'10','num')
    ,@(10,6,'NO','txt')
    ,@(11,1,'3164','forcetxt')
    ,@(11,6,'NO','txt')
    ,@(12,1,'3166','forcetxt')
    ,@(12,2,'10','num')
    ,@(12,3,'0','forcetxt')
    ,@(12,4,'0','forcetxt')
    ,@(12,5,'9.09%','forcetxt')
    ,@(12,6,'NO','txt')
    ,@(13,1,'3183','forcetxt')
    ,@(13,6,'NO','txt')
    ,@(14,1,'3309','forcetxt')
    ,@(14,2,'10','num')
    ,@(14,3,'0','forcetxt')
    ,@(14,4,'0','forcetxt')
    ,@(14,5,'0.67%','forcetxt')
    ,@(14,6,'NO','txt')
    ,@(15,1,'3310','forcetxt')
    ,@(15,6,'NO','txt')
    ,@(16,1,'3482','forcetxt')
    ,@(16,2,'8','num')
    ,@(16,6,'NO','txt')
    ,@(17,1,'3483','forcetxt')
    ,@(17,2,'8','num')
    ,@(17,6,'NO','txt')
    ,@(18,1,'3560','forcetxt')
    ,@(18,2,'8','num')
    ,@(18,6,'NO','txt')
    ,@(19,1,'3561','forcetxt')
    ,@(19,2,'9','num')
    ,@(19,6,'NO','txt')
    ,@(20,1,'3623','forcetxt')
    ,@(20,6,'NO','txt')
    ,@(21,1,'3707','forcetxt')
    ,@(21,2,'10','num')
    ,@(21,3,'1','forcetxt')
    ,@(21,4,'0','forcetxt')
    ,@(21,5,'2.63%','forcetxt')
    ,@(21,6,'NO','txt')
    ,@(22,1,'3711','forcetxt')
    ,@(22,6,'NO','txt')
    ,@(23,1,'3721','forcetxt')
    ,@(23,2,'10','num')
    ,@(23,3,'0','forcetxt')
    ,@(23,4,'0','forcetxt')
    ,@(23,6,'NO','txt')
    ,@(24,1,'3722','forcetxt')
    ,@(24,6,'NO','txt')
    ,@(25,1,'3730','forcetxt')
    ,@(25,2,'10','num')
    ,@(25,3,'0','forcetxt')
    ,@(25,4,'0','forcetxt')
    ,@(25,6,'NO','txt')
    ,@(26,1,'3754','forcetxt')
    ,@(26,2,'10','num')
    ,@(26,3,'0','forcetxt')
    ,@(26,4,'0','forcetxt')
    ,@(26,6,'NO','txt')
    ,@(27,1,'3759','forcetxt')
    ,@(27,2,'10','num')
    ,@(27,3,'1','forcetxt')
    ,@(27,4,'0','forcetxt')
    ,@(27,5,'7.11%','forcetxt')
    ,@(27,6,'NO','txt')
    ,@(28,1,'3764','forcetxt')
    ,@(28,2,'10','num')
    ,@(28,3,'0','forcetxt')
    ,@(28,4,'0','forcetxt')
    ,@(28,5,'4.93%','forcetxt')
    ,@(28,6,'NO','txt')
    ,@(29,1,'3773','forcetxt')
    ,@(29,6,'NO','txt')
    ,@(30,1,'3778','forcetxt')
    ,@(30,2,'9','num')
    ,@(30,3,'0','forcetxt')
    ,@(30,4,'0','forcetxt')
    ,@(30,6,'NO','txt')
    ,@(31,1,'3785','forcetxt')
    ,@(31,6,'NO','txt')
    ,@(32,1,'3912','forcetxt')
    ,@(32,6,'NO','txt')
    ,@(33,1,'3916','forcetxt')
    ,@(33,2,'11','num')
    ,@(33,3,'0','forcetxt')
    ,@(33,4,'0','forcetxt')
    ,@(33,5,'0.48%','forcetxt')
    ,@(33,6,'NO','txt')
    ,@(34,1,'4290','forcetxt')
    ,@(34,2,'10','num')
    ,@(34,3,'0','forcetxt')
    ,@(34,4,'1','forcetxt')
    ,@(34,5,'3.95%','forcetxt')
    ,@(34,6,'NO','txt')
    ,@(35,1,'4306','forcetxt')
    ,@(35,2,'11','num')
    ,@(35,3,'1','forcetxt')
    ,@(35,4,'0','forcetxt')
    ,@(35,5,'4.07%','forcetxt')
    ,@(35,6,'NO','txt')
    ,@(36,1,'4309','forcetxt')
    ,@(36,2,'10','num')
    ,@(36,3,'0','forcetxt')
    ,@(36,4,'0','forcetxt')
    ,@(36,6,'NO','txt')
    ,@(37,1,'4315','forcetxt')
    ,@(37,2,'10','num')
    ,@(37,3,'0','forcetxt')
    ,@(37,4,'0','forcetxt')
    ,@(37,5,'0.44%','forcetxt')
    ,@(37,6,'NO','txt')
    ,@(38,1,'4323','forcetxt')
    ,@(39,1,'4340','forcetxt')
)

$data5 = @(
    ,@(1,1,'MATCH_CODE','txt')
    ,@(1,2,'MAIDEN_OVERS','txt')
    ,@(1,3,'PERCENT_WICKETS_OF_ALL','txt')
    ,@(2,1,'2964','forcetxt')
    ,@(3,1,'3007','forcetxt')
    ,@(3,2,'0','forcetxt')
    ,@(3,3,'10.00%','forcetxt')
    ,@(4,1,'3008','forcetxt')
    ,@(5,1,'3135','forcetxt')
    ,@(6,1,'3138','forcetxt')
    ,@(6,2,'0','forcetxt')
    ,@(6,3,'30.00%','forcetxt')
    ,@(7,1,'3142','forcetxt')
    ,@(7,2,'1','forcetxt')
    ,@(7,3,'30.00%','forcetxt')
    ,@(8,1,'3147','forcetxt')
    ,@(9,1,'3149','forcetxt')
    ,@(10,1,'3153','forcetxt')
    ,@(10,2,'1','forcetxt')
    ,@(10,3,'10.00%','forcetxt')
    ,@(11,1,'3164','forcetxt')
    ,@(11,2,'0','forcetxt')
    ,@(11,3,'10.00%','forcetxt')
    ,@(12,1,'3166','forcetxt')
    ,@(13,1,'3183','forcetxt')
    ,@(13,2,'4','forcetxt')
    ,@(13,3,'40.00%','forcetxt')
    ,@(14,1,'3309','forcetxt')
    ,@(15,1,'3310','forcetxt')
    ,@(15,2,'0','forcetxt')
    ,@(15,3,'10.00%','forcetxt')
    ,@(16,1,'3482','forcetxt')
    ,@(17,1,'3483','forcetxt')
    ,@(17,2,'1','forcetxt')
    ,@(17,3,'10.00%','forcetxt')
    ,@(18,1,'3560','forcetxt')
    ,@(19,1,'3561','forcetxt')
    ,@(20,1,'3623','forcetxt')
    ,@(20,2,'1','forcetxt')
    ,@(20,3,'20.00%','forcetxt')
    ,@(21,1,'3707','forcetxt')
    ,@(22,1,'3711','forcetxt')
    ,@(22,2,'1','forcetxt')
    ,@(22,3,'50.00%','forcetxt')
    ,@(23,1,'3721','forcetxt')
    ,@(24,1,'3722','forcetxt')
    ,@(24,2,'0','forcetxt')
    ,@(24,3,'10.00%','forcetxt')
    ,@(25,1,'3730','forcetxt')
    ,@(26,1,'3754','forcetxt')
    ,@(26,2,'0','forcetxt')
    ,@(26,3,'20.00%','forcetxt')
    ,@(27,1,'3759','forcetxt')
    ,@(27,2,'0','forcetxt')
    ,@(27,3,'30.00%','forcetxt')
    ,@(28,1,'3764','forcetxt')
    ,@(28,2,'1','forcetxt')
    ,@(28,3,'10.00%','forcetxt')
    ,@(29,1,'3773','forcetxt')
    ,@(29,2,'0','forcetxt')
    ,@(29,3,'10.00%','forcetxt')
    ,@(30,1,'3778','forcetxt')
    ,@(31,1,'3785','forcetxt')
    ,@(31,2,'0','forcetxt')
    ,@(31,3,'10.00%','forcetxt')
    ,@(32,1,'3916','forcetxt')
    ,@(33,1,'4290','forcetxt')
    ,@(33,2,'1','forcetxt')
    ,@(33,3,'20.00%','forcetxt')
    ,@(34,1,'4306','forcetxt')
    ,@(34,2,'2','forcetxt')
    ,@(35,1,'4309','forcetxt')
    ,@(35,2,'0','forcetxt')
    ,@(35,3,'10.00%','forcetxt')
    ,@(36,1,'4315','forcetxt')
    ,@(36,2,'0','forcetxt')
    ,@(37,1,'4323','forcetxt')
    ,@(37,2,'1','forcetxt')
    ,@(38,1,'4340','forcetxt')
    ,@(38,2,'0','forcetxt')
)

function Write-SheetData($ws, $rows) {
    foreach ($row in $rows) {
        $r = $row[0]
        $c = $row[1]
        $val = $row[2]
        $kind = $row[3]
        $cell = $ws.Cells.Item($r, $c)
        if ($kind -eq "num") {
            $cell.Value = [double]$val
        } elseif ($kind -eq "forcetxt") {
            # Numeric-looking value that must stay text (MATCH_CODE / NUM_4 / NUM_6 / percentages)
            $cell.NumberFormat = "@"
            $cell.Value = $val
        } else {
            $cell.Value = $val
        }
    }
}

Write-SheetData $ws4 $data4
Write-SheetData $ws5 $data5

# Style the header rows like the rest of the workbook (bold, centered, bordered)
$hdr4 = $ws4.Range("A1:F1")
$hdr4.Font.Bold = $true
$hdr4.HorizontalAlignment = -4108
$hdr4.VerticalAlignment = -4160
$hdr4.Borders.LineStyle = 1

$hdr5 = $ws5.Range("A1:C1")
$hdr5.Font.Bold = $true
$hdr5.HorizontalAlignment = -4108
$hdr5.VerticalAlignment = -4160
$hdr5.Borders.LineStyle = 1

# Restore the originally active sheet/selection
$wb.Worksheets.Item("Player Info").Activate() | Out-Null
$wb.Worksheets.Item("Player Info").Range("A1").Select() | Out-Null
